# "text part of diploma done"
#
# The diploma's title-page paragraph 3 ("3. Вихідні дані до роботи: ...")
# previously ended with one long, single underlined run of spaces before the
# trailing "_______________________" blank.  The edit splits that run in two
# and drops Word's transient "_GoBack" (last-edit-position) bookmark right in
# the middle of it, moving it there from the empty paragraph above the
# signature block where it used to live.  Paragraph 4's enumeration run
# ". 7. " (right before "Висновки") is likewise split into ". 7." + " ".
#
# Because Word auto-renumbers every w:bookmarkStart/w:bookmarkEnd w:id in
# document order whenever the package is saved, we only need to get bookmark
# *names* and *positions* right -- the numeric ids fall out on their own.

$d = $word.ActiveDocument

# --- 1. Drop the old "_GoBack" bookmark (sat alone in an otherwise-empty
#        paragraph right above the signature block). ---
$d.Bookmarks.Item("_GoBack").Delete()

# --- 2. Split the long underlined whitespace run in the "Вихідні дані"
#        paragraph, 27 characters in ("." + 26 spaces | 28 spaces), and drop
#        a fresh (empty-range) "_GoBack" bookmark right at that split point.
#        Word/the OOXML writer splits the run into two <w:r> automatically
#        when a bookmark lands inside it. ---
$rWhitespace = $d.Content
$rWhitespace.Find.Execute(".                                                      ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $rWhitespace.Start + 27
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

# --- 3. Split ". 7. " (before "Висновки") into ". 7." and " " as two
#        separate runs with identical formatting.  There's no bookmark left
#        behind here, so add a throwaway one purely to force the run split,
#        then delete it again. ---
$rEnum = $d.Content
$rEnum.Find.Execute(". 7. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $rEnum.Start + 4
$d.Bookmarks.Add("ZZZ_TempRunSplit", $d.Range($splitPos, $splitPos))
$d.Bookmarks.Item("ZZZ_TempRunSplit").Delete()
